$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
Write-Host "Default numberformat sample: $($ws.Cells.Item(5,5).NumberFormat)"
$ws.Cells.Item(1000,26).NumberFormat = $ws.Cells.Item(5,5).NumberFormat
